# Add the "subteam" project data to the Electrical sheet and make it the
# active tab (addition of the subteam page, still in development).

$wb = $excel.ActiveWorkbook
$projects = $wb.Worksheets.Item(1)
$electrical = $wb.Worksheets.Item(2)

# --- Fill in the four placeholder "subteam" projects on the Electrical sheet ---
# Enter column-by-column (Title, then Blurb, then Big Blurb) so new shared
# strings land in the same order the workbook author typed them in.
$electrical.Range("A2").Value = "Project 1"
$electrical.Range("A3").Value = "Project 2"
$electrical.Range("A4").Value = "Project 3"
$electrical.Range("A5").Value = "Project 4"

$electrical.Range("C2").Value = "This is project 1"
$electrical.Range("C3").Value = "This is project 2"
$electrical.Range("C4").Value = "This is project 3"
$electrical.Range("C5").Value = "This is project 4"

$electrical.Range("D2").Value = "big blurb for project 1"
$electrical.Range("D3").Value = "big blurb for project 2"
$electrical.Range("D4").Value = "big blurb for project 3"
$electrical.Range("D5").Value = "big blurb for project 4"

# --- Switch focus onto the new Electrical (subteam) page ---
$electrical.Activate() | Out-Null
$electrical.Range("D2:D5").Select() | Out-Null

Write-Output "Applied subteam page edits"
